$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-cell formatting (style index 2: bold, bordered, centered, date numfmt)
# from the last existing row (A328) so the new date cells inherit the same style.
$ws.Range("A328").Copy()
$ws.Range("A329:A343").PasteSpecial(-4122)

$ws.Range("A329").Value = 44403
$ws.Range("B329").Value = 4
$ws.Range("C329").Value = 25
$ws.Range("D329").Value = 72.74631903625676

$ws.Range("A330").Value = 44404
$ws.Range("B330").Value = 2
$ws.Range("C330").Value = 27
$ws.Range("D330").Value = 78.56602455915731

$ws.Range("A331").Value = 44405
$ws.Range("B331").Value = 6
$ws.Range("C331").Value = 31
$ws.Range("D331").Value = 90.2054356049584

$ws.Range("A332").Value = 44406
$ws.Range("B332").Value = 9
$ws.Range("C332").Value = 34
$ws.Range("D332").Value = 98.9349938893092

$ws.Range("A333").Value = 44407
$ws.Range("B333").Value = 5
$ws.Range("C333").Value = 37
$ws.Range("D333").Value = 107.66455217366

$ws.Range("A334").Value = 44408
$ws.Range("B334").Value = 3
$ws.Range("C334").Value = 34
$ws.Range("D334").Value = 98.9349938893092

$ws.Range("A335").Value = 44409
$ws.Range("B335").Value = 0
$ws.Range("C335").Value = 29
$ws.Range("D335").Value = 84.38573008205785

$ws.Range("A336").Value = 44410
$ws.Range("B336").Value = 5
$ws.Range("C336").Value = 30
$ws.Range("D336").Value = 87.29558284350811

$ws.Range("A337").Value = 44411
$ws.Range("B337").Value = 2
$ws.Range("C337").Value = 30
$ws.Range("D337").Value = 87.29558284350811

$ws.Range("A338").Value = 44412
$ws.Range("B338").Value = 1
$ws.Range("C338").Value = 25
$ws.Range("D338").Value = 72.74631903625676

$ws.Range("A339").Value = 44413
$ws.Range("B339").Value = 4
$ws.Range("C339").Value = 20
$ws.Range("D339").Value = 58.19705522900541

$ws.Range("A340").Value = 44414
$ws.Range("B340").Value = 5
$ws.Range("C340").Value = 20
$ws.Range("D340").Value = 58.19705522900541

$ws.Range("A341").Value = 44415
$ws.Range("B341").Value = 9
$ws.Range("C341").Value = 26
$ws.Range("D341").Value = 75.65617179770703

$ws.Range("A342").Value = 44416
$ws.Range("B342").Value = 10
$ws.Range("C342").Value = 36
$ws.Range("D342").Value = 104.7546994122097

$ws.Range("A343").Value = 44417
$ws.Range("B343").Value = 4
$ws.Range("C343").Value = 35
$ws.Range("D343").Value = 101.8448466507595
